$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-16 (A=8..14) get new scheme labels: Gaussian-Quadrature and the three
# new Spiral-* schemes are inserted before NoRotation-tilt60deg, pushing the
# remaining rotation/hex-grid schemes down.
$ws.Range("B10").Value2 = "Gaussian-Quadrature"
$ws.Range("B11").Value2 = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value2 = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value2 = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value2 = "NoRotation-tilt60deg"
$ws.Range("B15").Value2 = "Rotation-NoTilt"
$ws.Range("B16").Value2 = "Rotation-60detTilt"

# New rows 17-19 (A=15..17) carry the schemes that got pushed past the old
# bottom of the table.
$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "HexGrid-90degTilt5degRes"
$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = "HexGrid-60degTilt5degRes"

# Fill in the averaged-intensity values (all 1) for the 14 HKL columns on the
# three new rows, matching every other data row.
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P")
foreach ($r in 17..19) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = 1
    }
}

# Copy column-A number formatting (bold, centered, bordered) onto the three
# new label cells so they match the rest of the table.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
